# 自动更新Excel文件 - 2025-11-05 23:12:38
# Decrement the "剩余" (remaining) count in column E by 1 for every data row
# (rows 2-99) to reflect one more day elapsed, except:
#   - row 36, whose start date (column F) is malformed/unparsed and whose
#     remaining count already equals the total, so it is left untouched
#   - row 94 ("刘记水饺"), whose subscription was renewed: remaining (E) is
#     reset to the total (7) and the start date (F) is bumped to 20251106

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }

    if ($r -eq 94) {
        $ws.Cells.Item($r, 5).Value = 7
        $ws.Cells.Item($r, 6).Value = 20251106
        continue
    }

    $remainingCell = $ws.Cells.Item($r, 5)
    $current = $remainingCell.Value2
    $remainingCell.Value = $current - 1
}
